$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.289.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6252"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07476"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07711"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.840.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.981"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6767"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001025"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.100.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.090"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.329.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.369"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1378"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.365"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.400"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05704"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.094"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.818"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.141"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6912"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.586"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.240.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.503"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9040"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9981"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.001.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.074"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1167"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.976"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3934"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000113"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.51%  "
